$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two early test-case rows (old rows 2 and 3), shifting remaining
# rows up. This leaves row 1 (header/template row) intact and turns the
# former rows 4-8 into rows 2-6.
$ws.Rows("2:3").Delete()
